$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly generated data points (value 0) for previously blank rows.
$ws.Range("A24").Value = 0
$ws.Range("A41").Value = 0
$ws.Range("A43").Value = 0
$ws.Range("A54").Value = 0
$ws.Range("A55").Value = 0
$ws.Range("A56").Value = 0

# Restore the view state: scrolled near the bottom of the data with A56 selected.
$excel.ActiveWindow.ScrollRow = 388
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A56").Select()
